$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'capri yoga pants for women gray'
$ws.Cells.Item(2, 1).Value = 'capri yoga pants for women green'
$ws.Cells.Item(3, 1).Value = 'capri yoga pants for women high waist'
$ws.Cells.Item(4, 1).Value = 'capri yoga pants xl'
$ws.Cells.Item(5, 1).Value = 'capri yoga tights women'
$ws.Cells.Item(6, 1).Value = 'caprice leggings'
$ws.Cells.Item(7, 1).Value = 'capris 16w'
$ws.Cells.Item(8, 1).Value = 'capris 2x'
$ws.Cells.Item(9, 1).Value = 'capris 5x'
$ws.Cells.Item(10, 1).Value = 'capris and shorts'
$ws.Cells.Item(11, 1).Value = 'capris athletic'
$ws.Cells.Item(12, 1).Value = 'capris black'
$ws.Cells.Item(13, 1).Value = 'capris clearance'
$ws.Cells.Item(14, 1).Value = 'capris for short women'
$ws.Cells.Item(15, 1).Value = 'capris for women'
$ws.Cells.Item(16, 1).Value = 'capris for women high waist'
$ws.Cells.Item(17, 1).Value = 'capris for women leggings'
$ws.Cells.Item(18, 1).Value = 'capris for women plus size'
$ws.Cells.Item(19, 1).Value = 'capris for women size 20'
$ws.Cells.Item(20, 1).Value = 'capris high waist'
$ws.Cells.Item(21, 1).Value = 'capris ladies'
$ws.Cells.Item(22, 1).Value = 'capris leggings'
$ws.Cells.Item(23, 1).Value = 'capris leggings for women'
$ws.Cells.Item(24, 1).Value = 'capris leggings with pockets'
$ws.Cells.Item(25, 1).Value = 'capris pants'
$ws.Cells.Item(26, 1).Value = 'capris pants for women'
$ws.Cells.Item(27, 1).Value = 'capris shorts'
$ws.Cells.Item(28, 1).Value = 'capris shorts for women'
$ws.Cells.Item(29, 1).Value = 'capris size 2'
$ws.Cells.Item(30, 1).Value = 'capris size 4'
$ws.Cells.Item(31, 1).Value = 'capris skirt'
$ws.Cells.Item(32, 1).Value = 'capris stretch'
$ws.Cells.Item(33, 1).Value = 'capris tall'
$ws.Cells.Item(34, 1).Value = 'capris tights'
$ws.Cells.Item(35, 1).Value = 'capris tights for women'
$ws.Cells.Item(36, 1).Value = 'capris with holes'
$ws.Cells.Item(37, 1).Value = 'capris with holes for women'
$ws.Cells.Item(38, 1).Value = 'capris with skirt'
$ws.Cells.Item(39, 1).Value = 'capris workout leggings for women'
$ws.Cells.Item(40, 1).Value = 'capris workout pants'
$ws.Cells.Item(41, 1).Value = 'capris workout pants women'
$ws.Cells.Item(42, 1).Value = 'capris yoga'
$ws.Cells.Item(43, 1).Value = 'capris yoga pants'
$ws.Cells.Item(44, 1).Value = 'capris yoga pants for women'
$ws.Cells.Item(45, 1).Value = 'car leggings'
$ws.Cells.Item(46, 1).Value = 'car recovery gear'
$ws.Cells.Item(47, 1).Value = 'careless body wash'
$ws.Cells.Item(48, 1).Value = 'carhartt force utility knit legging'
$ws.Cells.Item(49, 1).Value = 'cars training pants'
$ws.Cells.Item(50, 1).Value = 'cashmere travel sets'
$ws.Cells.Item(51, 1).Value = 'casual black leggings for women'
$ws.Cells.Item(52, 1).Value = 'casual woman tights'
$ws.Cells.Item(53, 1).Value = 'cellulite compression leggings'
$ws.Cells.Item(54, 1).Value = 'cellulite compression leggings women'
$ws.Cells.Item(55, 1).Value = 'cellulite leggings for women compression'
$ws.Cells.Item(56, 1).Value = 'cellulite on thighs'
$ws.Cells.Item(57, 1).Value = 'cellulite thighs'
$ws.Cells.Item(58, 1).Value = 'cep women'
$ws.Cells.Item(59, 1).Value = 'chaffing balls'
$ws.Cells.Item(60, 1).Value = 'chaffing stick'
$ws.Cells.Item(61, 1).Value = 'chafing bands for thighs'
$ws.Cells.Item(62, 1).Value = 'chafing runners'
$ws.Cells.Item(63, 1).Value = 'chafing running'
$ws.Cells.Item(64, 1).Value = 'champion compression pants women'
$ws.Cells.Item(65, 1).Value = 'cheap black leggings for women'
$ws.Cells.Item(66, 1).Value = 'cheap capri leggings'
$ws.Cells.Item(67, 1).Value = 'cheap green tights'
$ws.Cells.Item(68, 1).Value = 'cheap lululemon leggings'
$ws.Cells.Item(69, 1).Value = 'cheap tights'
$ws.Cells.Item(70, 1).Value = 'cheap tights for women'
$ws.Cells.Item(71, 1).Value = 'cheap trick apparel'
$ws.Cells.Item(72, 1).Value = 'cheek art'
$ws.Cells.Item(73, 1).Value = 'cheeks fit body'
$ws.Cells.Item(74, 1).Value = 'cheerleader leggings'
$ws.Cells.Item(75, 1).Value = 'cheerleader shorts for women'
$ws.Cells.Item(76, 1).Value = 'cheerleader tights'
$ws.Cells.Item(77, 1).Value = 'cheerleader tights women'
$ws.Cells.Item(78, 1).Value = 'cheerleader training'
$ws.Cells.Item(79, 1).Value = 'chiropractors blend'
$ws.Cells.Item(80, 1).Value = 'chocolate basketballs'
$ws.Cells.Item(81, 1).Value = 'choice apparel basketball shorts'
$ws.Cells.Item(82, 1).Value = 'chronic pain your key to recovery'
$ws.Cells.Item(83, 1).Value = 'cigarette bomb'
$ws.Cells.Item(84, 1).Value = 'cigarette pant'
$ws.Cells.Item(85, 1).Value = 'cigarette pants women'
$ws.Cells.Item(86, 1).Value = 'cigarette post'
$ws.Cells.Item(87, 1).Value = 'cigarette trousers'
$ws.Cells.Item(88, 1).Value = 'cigarettes chocolate'
$ws.Cells.Item(89, 1).Value = 'circulation leggings women'
$ws.Cells.Item(90, 1).Value = 'circulation thighs'
$ws.Cells.Item(91, 1).Value = 'circulation tights'
$ws.Cells.Item(92, 1).Value = 'class rings'
$ws.Cells.Item(93, 1).Value = 'class rings for women'
$ws.Cells.Item(94, 1).Value = 'class struggle game'
$ws.Cells.Item(95, 1).Value = 'clear leggings for women'
$ws.Cells.Item(96, 1).Value = 'climbing injury free'
$ws.Cells.Item(97, 1).Value = 'climbing oants'
$ws.Cells.Item(98, 1).Value = 'climbing skin glue'
$ws.Cells.Item(99, 1).Value = 'climbing training balls'
$ws.Cells.Item(100, 1).Value = 'close fitted dresses for women'
